$d = $word.ActiveDocument

$replacements = @(
    @{old = "891÷3=297, 0"; new = "170÷4=42, 2"},
    @{old = "983÷4=245, 3"; new = "381÷3=127, 0"},
    @{old = "637÷2=318, 1"; new = "602÷2=301, 0"},
    @{old = "948÷3=316, 0"; new = "128÷3=42, 2"},
    @{old = "426÷9=47, 3"; new = "339÷8=42, 3"},
    @{old = "834÷2=417, 0"; new = "559÷5=111, 4"},
    @{old = "227÷9=25, 2"; new = "916÷2=458, 0"},
    @{old = "332÷9=36, 8"; new = "268÷7=38, 2"},
    @{old = "678÷6=113, 0"; new = "718÷2=359, 0"},
    @{old = "589÷3=196, 1"; new = "305÷3=101, 2"},
    @{old = "197÷8=24, 5"; new = "320÷8=40, 0"},
    @{old = "148÷9=16, 4"; new = "115÷4=28, 3"},
    @{old = "249÷5=49, 4"; new = "314÷8=39, 2"},
    @{old = "169÷6=28, 1"; new = "715÷7=102, 1"},
    @{old = "377÷6=62, 5"; new = "771÷6=128, 3"},
    @{old = "786÷5=157, 1"; new = "563÷4=140, 3"},
    @{old = "731÷5=146, 1"; new = "338÷4=84, 2"},
    @{old = "974÷9=108, 2"; new = "924÷4=231, 0"},
    @{old = "359÷6=59, 5"; new = "343÷4=85, 3"},
    @{old = "669÷7=95, 4"; new = "860÷6=143, 2"},
    @{old = "459÷5=91, 4"; new = "362÷8=45, 2"},
    @{old = "698÷7=99, 5"; new = "189÷2=94, 1"},
    @{old = "919÷3=306, 1"; new = "308÷2=154, 0"},
    @{old = "517÷3=172, 1"; new = "956÷6=159, 2"},
    @{old = "405÷3=135, 0"; new = "520÷8=65, 0"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
